# Update "南宁-漫展信息" workbook: a new exhibition entry
# ("南宁·THOsp·幻想朱槿绘翠5") is inserted as row 3 in both the "展览"
# and "全部类型" sheets, pushing the following rows down by one and
# refreshing a handful of visitor-count / min-price values along the way.

$wb = $excel.ActiveWorkbook

function Set-TextValue($ws, $cellRef, $text) {
    # Plain `.Value = $text` lets Excel auto-detect date-like strings
    # ("2024-05-03") and coerce them into real dates, which also stamps a
    # new number-format style onto the cell. Force the format to Text
    # first so the literal string is kept, then paste a known-plain
    # (style-0, General-format) cell's formatting back over it so the
    # cell's style index ends up untouched/default, same as every other
    # plain-text cell in the sheet.
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $text
    $ws.Range("B2").Copy() | Out-Null
    $range.PasteSpecial(-4122) | Out-Null
}

function Update-Sheet($sheetName, $lastRow) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Insert a fresh row 3; everything from the old row 3 onward shifts
    # down by one row (row4->5, row5->6, ...).
    $ws.Rows.Item(3).Insert()

    # The insert carries the row-above's formatting onto column A of the
    # new row (borders/bold/alignment) but Excel represents it with a
    # different style index than the rest of column A; normalise it back
    # to match.
    $ws.Range("A2").Copy() | Out-Null
    $ws.Range("A3").PasteSpecial(-4122) | Out-Null

    # F2 (想去人数 for the first / headline listing) ticks up.
    $ws.Range("F2").Value = 5622

    # Row 3: brand-new listing.
    Set-TextValue $ws "B3" "2024-05-03"
    $ws.Range("C3").Value = "南宁·THOsp·幻想朱槿绘翠5"
    $ws.Range("D3").Value = "金湖路58号 广西建设五象大酒店"
    $ws.Range("E3").Value = "2024.05.03 09:00-05.03 22:00"
    $ws.Range("F3").Value = 1
    $ws.Range("G3").Value = 60
    $ws.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=84483"
    $ws.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202404/HIiFmh7F1713270854919.png"

    # Row 4: previously row 3 (南宁·原x穹x崩only), F bumps 171 -> 172.
    Set-TextValue $ws "B4" "2024-05-19"
    $ws.Range("C4").Value = "南宁·原x穹x崩only"
    $ws.Range("D4").Value = "明秀东路157号 利泰国际大酒店"
    $ws.Range("E4").Value = "2024.05.19 10:00-05.19 17:00"
    $ws.Range("F4").Value = 172
    $ws.Range("G4").Value = 35
    $ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=83070"
    $ws.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202403/I8tScigE1710918412731.jpeg"

    # Row 5: previously row 4 (南宁·布谷鸟动漫展4th), F bumps 951 -> 956.
    Set-TextValue $ws "B5" "2024-06-09"
    $ws.Range("C5").Value = "南宁·布谷鸟动漫展4th"
    $ws.Range("D5").Value = "亭洪路45号 百益上河城"
    $ws.Range("E5").Value = "2024.06.09 10:00-06.10 17:00"
    $ws.Range("F5").Value = 956
    $ws.Range("G5").Value = 50
    $ws.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=82241"
    $ws.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202403/uzZqZov91709281147333.jpeg"

    # Row 6: previously row 5 (南宁·恋与深空only), F bumps 18 -> 20.
    Set-TextValue $ws "B6" "2024-06-09"
    $ws.Range("C6").Value = "南宁·恋与深空only"
    $ws.Range("D6").Value = "新阳路227号南宁第三人民医院旁新秀佳园对面 卡尔顿东方银龙酒店"
    $ws.Range("E6").Value = "2024.06.09 10:00-06.09 17:00"
    $ws.Range("F6").Value = 20
    $ws.Range("G6").Value = 50
    $ws.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=84444"
    $ws.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202404/6ZVHU1F91713340880421.jpeg"

    # Re-number the serial index column (A) for every data row so it
    # stays 1,2,3,... after the insert (inserted rows keep whatever value
    # shifted into them, which is wrong for row 3 and everything below).
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Range("A$r").Value = $r - 1
    }
}

# "展览" originally had 5 data rows (A1:I5) -> 6 after the insert.
Update-Sheet "展览" 6

# "全部类型" originally had 6 data rows (A1:I6) -> 7 after the insert;
# its former row 6 (南宁·浪漫古典·百年经典世界名曲音乐会) simply slides
# down to row 7 untouched by the insert itself.
Update-Sheet "全部类型" 7
